$wb = $excel.ActiveWorkbook

$wsSalesman = $wb.Worksheets.Item("Salesman_Add")
$wsSalesman.Columns.Item(2).Delete()
$wsSalesman.Range("A1").Value = "Employee Code/UserId"

$wsSalesman.Activate()
